$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.896.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.299.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.511"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.660.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.306.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.783"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.837.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0689"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.006.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.34%  "
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.529.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "
